# Restructure ontology: remove mfd_hab1=Urban if mfd_areatype=Urban
#
# Rows 10-17 are the biofilm / Sediment / mfd_areatype=Urban samples. For
# these rows:
#   - habitat_typenumber (F) is bumped to the next EUNIS-ish code
#     (2130 -> 2300, 2120 -> 2200).
#   - mfd_hab1 (N) should no longer just repeat mfd_areatype's "Urban"
#     value, so it takes on what used to be mfd_hab2's value ("Other").
#   - mfd_hab2 (O) takes on what used to be mfd_hab3's value (the specific
#     habitat note, e.g. "High chalk concentration (limestone quarry)").
#   - mfd_hab3 (P) is removed entirely (the row now only has hab1/hab2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New habitat_typenumber (column F) values, keyed by row.
$newHabitatTypeNumber = @{
    10 = "2300"
    11 = "2300"
    12 = "2300"
    13 = "2300"
    14 = "2200"
    15 = "2300"
    16 = "2300"
    17 = "2300"
}

for ($row = 10; $row -le 17; $row++) {
    # F: habitat_typenumber gets bumped. Prefix with an apostrophe so Excel
    # keeps storing this numeric-looking value as text, matching the
    # original cell's (inline string) data type.
    $ws.Range("F$row").Value = "'" + $newHabitatTypeNumber[$row]

    # Shift mfd_hab2 -> mfd_hab1 and mfd_hab3 -> mfd_hab2, then clear
    # mfd_hab3 so the column disappears from that row.
    $hab2 = $ws.Range("O$row").Value2
    $hab3 = $ws.Range("P$row").Value2

    $ws.Range("N$row").Value = $hab2
    $ws.Range("O$row").Value = $hab3
    $ws.Range("P$row").Value = $null
}
